$d = $word.ActiveDocument

# Locate the paragraph "a = toggle auto-naming neurons".
$targetIndex = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "toggle auto-naming neurons") {
        $targetIndex = $i
        break
    }
}

$target = $d.Paragraphs.Item($targetIndex)

# Duplicate the whole paragraph (runs + formatting, via copy/paste) so the
# new paragraph inherits identical run/paragraph formatting (Arial, bold
# single-letter shortcut run, line spacing 360/auto) without having to
# hand-build run properties on freshly inserted text.
$srcRange = $d.Range($target.Range.Start, $target.Range.End)
$srcRange.Copy()
$pasteAt = $d.Range($target.Range.End, $target.Range.End)
$pasteAt.Paste()

$newPara = $d.Paragraphs.Item($targetIndex + 1)

# Rewrite the bold shortcut-key run: "a" -> "b".
$letterRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)
$letterRange.Text = "b"

# Rewrite the remaining (non-bold) run's text.
$restRange = $d.Range($newPara.Range.Start + 1, $newPara.Range.End - 1)
$restRange.Text = " = toggle showing neuron birth times"
